$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted before the existing row 369,
# pushing all subsequent rows (old 369..474) down by one (new 370..475).
$ws.Rows.Item(369).Insert()

# Populate the newly inserted row 369 with the latest week's data.
$ws.Cells.Item(369, 1).Value = 4
$ws.Cells.Item(369, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(369, 3).Value = "Los Lagos"
$ws.Cells.Item(369, 4).Value = 45093
$ws.Cells.Item(369, 5).Value = 10
$ws.Cells.Item(369, 6).Value = 100112003
$ws.Cells.Item(369, 7).Value = "Ajo"
$ws.Cells.Item(369, 8).Value = "Chino"
$ws.Cells.Item(369, 9).Value = "Primera"
$ws.Cells.Item(369, 10).Value = 240
$ws.Cells.Item(369, 11).Value = 19000
$ws.Cells.Item(369, 12).Value = 21000
$ws.Cells.Item(369, 13).Value = 20000
$ws.Cells.Item(369, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(369, 15).Value = "China"
$ws.Cells.Item(369, 16).Value = 2000
$ws.Cells.Item(369, 17).Value = 10
$ws.Cells.Item(369, 18).Value = "Hortaliza"
